$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" ---
$metadata = $wb.Worksheets.Item("Metadata")

# Version bump 5.0.0 -> 6.0.0
$metadata.Range("B3").Value = "6.0.0"

# Date bump
$metadata.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank; now populated
$metadata.Range("B9").Value = "Alvearie Team"

# Row 10 (was a duplicated "Contact" row) becomes "Jurisdiction"/"United States of America"
$metadata.Range("A10").Value = "Jurisdiction"
$metadata.Range("B10").Value = "United States of America"

# Row 11 was an exact duplicate ("Contact" / "No display for ContactDetail") - remove it entirely
$metadata.Rows.Item(11).Delete()

# --- Sheet 2: "Elements" ---
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short/Definition customized for this profile
$elements.Range("K2").Value = "Employee FTE Count"
$elements.Range("L2").Value = "Full Time Equivalent calculation for the employee. Example 1: A full-time employee would have a value of 1. Example 2: A 20-hour per week employee would have a value of 0.5."
